# Lab sign-in sheet: append sign-in rows for the lab sessions that were
# worked (4/22 Needle Demo + 4/24 Who are Doctors), matching the rows the
# LabSignup app writes to SignInSheet.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=FirstName B=LastName C=Title D=LabName E=LabDay F=LabStart
#          G=LabEnd H=LabSignInTime I=LabHours

$data = @(
    @("test",   "test",  "DENTAL",      "4/24/2023 - Who are Doctors", "4/24/2023 ", "4:00PM", "6:00PM", "4/28/2023 10:39:29 PM", "02:00:00"),
    @("justin",  "k",      "DENTAL",      "4/24/2023 - Who are Doctors", "4/24/2023 ", "4:00PM", "6:00PM", "4/28/2023 10:49:16 PM", "02:00:00"),
    @("tes",     "test",  "OCC THERAPY", "4/24/2023 - Who are Doctors", "4/24/2023 ", "4:00PM", "6:00PM", "4/28/2023 11:14:24 PM", "02:00:00"),
    @("test",    "test",  "DENTAL",      "4/22/2023 - Needle Demo",     "4/22/2023 ", "4:00PM", "6:45PM", "4/28/2023 11:32:48 PM", "02:45:00")
)

$row = 2
foreach ($rec in $data) {
    for ($col = 0; $col -lt $rec.Length; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $rec[$col]
    }
    $row++
}

# Touch the header/footer so the sheet carries an (empty) headerFooter part,
# same as what the app's own writer emits on save.
$ws.PageSetup.CenterHeader = ""
